# Update benchmark results data (Insertion Sort, Quicksort, Heap Sort rows)
# with refreshed timing numbers from the report script re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.887
$ws.Range("C2").Value = 35.3
$ws.Range("D2").Value = 118.405
$ws.Range("E2").Value = 256.676
$ws.Range("F2").Value = 530.1849999999999
$ws.Range("G2").Value = 791.255
$ws.Range("H2").Value = 3349.476
$ws.Range("I2").Value = 6866.516
$ws.Range("J2").Value = 11918.365
$ws.Range("K2").Value = 18299.652
$ws.Range("L2").Value = 26742.385
$ws.Range("M2").Value = 35897.407
$ws.Range("N2").Value = 48638.799

$ws.Range("B3").Value = 1.4
$ws.Range("C3").Value = 3.491
$ws.Range("D3").Value = 7.187
$ws.Range("E3").Value = 10.965
$ws.Range("F3").Value = 14.962
$ws.Range("G3").Value = 21.052
$ws.Range("H3").Value = 41.845
$ws.Range("I3").Value = 61.537
$ws.Range("J3").Value = 79.682
$ws.Range("K3").Value = 114.689
$ws.Range("L3").Value = 123.069
$ws.Range("M3").Value = 145.967
$ws.Range("N3").Value = 171.44

$ws.Range("B4").Value = 2.497
$ws.Range("C4").Value = 8.875999999999999
$ws.Range("D4").Value = 20.245
$ws.Range("E4").Value = 37.396
$ws.Range("F4").Value = 58.55
$ws.Range("G4").Value = 63.83
$ws.Range("H4").Value = 144.878
$ws.Range("I4").Value = 218.415
$ws.Range("J4").Value = 297.171
$ws.Range("K4").Value = 383.933
$ws.Range("L4").Value = 519.607
$ws.Range("M4").Value = 624.441
$ws.Range("N4").Value = 680.198
